# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# mirroring the rows produced by the automated SeniorConnect logger run on
# 2026-02-06.

$wb = $excel.ActiveWorkbook

function Set-LogCellValue {
    param(
        $Cell,
        $Text
    )

    # The source log stores every column as plain text (e.g. "2026-02-06",
    # "84.4%"). Excel's normal typing heuristics would otherwise turn a
    # "YYYY-MM-DD" value into a date serial, or a "NN.N%" value into a
    # fractional percentage number. Prefixing those look-alike values with
    # a quote keeps them as literal text, matching the rest of the log.
    if ($Text -match "^\d{4}-\d{2}-\d{2}$" -or $Text -match "%$") {
        $Cell.Value = "'" + $Text
    } else {
        $Cell.Value = $Text
    }
}

function Add-LogRows {
    param(
        $Worksheet,
        $StartRow,
        $Rows
    )

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $rowData = $Rows[$i]

        for ($col = 1; $col -le $rowData.Count; $col++) {
            Set-LogCellValue $Worksheet.Cells.Item($r, $col) $rowData[$col - 1]
        }
    }
}

# ---------------------------------------------------------------------
# PIR sheet: rows 86-98 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-06","09:47:09","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:10","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:15","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:20","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:25","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:30","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:35","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:40","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:45","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:50","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:47:55","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:48:00","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:48:05","09:00","Bathroom","No Motion","Inactive")
)
Add-LogRows $pir 86 $pirRows

# ---------------------------------------------------------------------
# Humidity sheet: rows 26-35
# ---------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-06","09:47:09","09:00","Bathroom","84.4%","Active"),
    @("2026-02-06","09:47:09","09:00","Bathroom","79.2%","Active"),
    @("2026-02-06","09:47:14","09:00","Bathroom","76.6%","Active"),
    @("2026-02-06","09:47:19","09:00","Bathroom","75.4%","Active"),
    @("2026-02-06","09:47:25","09:00","Bathroom","74.5%","Active"),
    @("2026-02-06","09:47:40","09:00","Bathroom","72.0%","Active"),
    @("2026-02-06","09:47:50","09:00","Bathroom","71.6%","Active"),
    @("2026-02-06","09:47:55","09:00","Bathroom","72.3%","Active"),
    @("2026-02-06","09:48:00","09:00","Bathroom","72.4%","Active"),
    @("2026-02-06","09:48:05","09:00","Bathroom","72.5%","Active")
)
Add-LogRows $humidity 26 $humidityRows

# ---------------------------------------------------------------------
# Temperature sheet: rows 26-35
# ---------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-06","09:47:09","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:47:10","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:47:15","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:47:20","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:47:25","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:47:40","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:47:50","09:00","Bathroom","27.4C","Active"),
    @("2026-02-06","09:47:55","09:00","Bathroom","27.4C","Active"),
    @("2026-02-06","09:48:00","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:48:05","09:00","Bathroom","27.5C","Active")
)
Add-LogRows $temperature 26 $temperatureRows
